# Auto-applies the RKI-COVID-19_Todesfaelle.xlsx update (downloaded 2020-12-21--13-35-01)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 1: COVID_Todesfaelle (per Sterbewoche / age-group index) ---
$ws1.Cells.Item(1, 1).Value = "Sterbewoche"
$ws1.Cells.Item(1, 2).Value = "Anzahl verstorbene COVID-19 Fälle"

$s1A = @(
1, 2, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47
)
$s1B = @(
"0", "<4", "0", "0", "<4", "<4", "0", "0", "0", "0", "<4", "18", "162", "601", "1369", "1740", "1594", "1168", "782", "514", "351", "271", "150", "112", "73", "47", "51", "46", "26", "26", "31", "30", "29", "32", "39", "37", "19", "32", "53", "65", "79", "116", "230", "387", "746", "1134", "1497", "1887"
)

for ($i = 0; $i -lt $s1A.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 1).Value = $s1A[$i]
    $bcell = $ws1.Cells.Item($row, 2)
    $bcell.NumberFormat = "@"
    $bcell.Value = $s1B[$i]
}

# --- Sheet 2: COVID_Todesfaelle_Monat (per SterbeMonat) ---
$ws2.Cells.Item(1, 1).Value = "SterbeMonat"
$ws2.Cells.Item(1, 2).Value = "Anzahl verstorbene COVID-19 Fälle"

$s2A = @(
3, 4, 5, 6, 7, 8, 9, 10, 11
)
$s2B = @(
"1119", "6041", "1562", "302", "131", "148", "201", "1398", "4644"
)

for ($i = 0; $i -lt $s2A.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $s2A[$i]
    $bcell = $ws2.Cells.Item($row, 2)
    $bcell.NumberFormat = "@"
    $bcell.Value = $s2B[$i]
}
